# Generate Report for handback
# Updates the zh-cn and de-de language sheets of the localization-status
# workbook: marks rows 2 & 3 as "Handed back", fills in the new
# "Latest Target File" / "Latest Handback File" columns (E/F) with
# hyperlinks, and stamps the "Latest Handback DateTime" column (G).
# The "Status" text is a shared string reused by the Overview roll-up
# sheet too, so that sheet's B2:C3 cells flip to the new wording as well.

$wb = $excel.ActiveWorkbook
$newStatus = "Handed back: in sync with en-US"

# Per-language handback timestamp (column G, rows 2 & 3).
$handbackInfo = @{
    "zh-cn" = "2016-01-19 06:47:45"
    "de-de" = "2016-01-19 06:48:03"
}

# The Overview sheet mirrors the per-language "Status" column (B = zh-cn,
# C = de-de) for rows 2 & 3; it shares the same underlying string so it
# picks up the new wording too.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in @(2, 3)) {
    $wsOverview.Cells.Item($row, 2).Value = $newStatus
    $wsOverview.Cells.Item($row, 3).Value = $newStatus
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Look up the existing hyperlink target/display for A2 ("Source File
    # Name") and C2 ("Latest Handoff File") so the new "Latest Target
    # File" (E) / "Latest Handback File" (F) columns can reuse them.
    $aAddr = $ws.Range("A2").Address()
    $cAddr = $ws.Range("C2").Address()
    $sourceUrl = ""
    $sourceDisplay = ""
    $handoffUrl = ""
    $handoffDisplay = ""
    foreach ($hl in $ws.Hyperlinks) {
        $hlAddr = $hl.Range.Address()
        if ($hlAddr -eq $aAddr) {
            $sourceUrl = $hl.Address
            $sourceDisplay = $hl.TextToDisplay
        }
        if ($hlAddr -eq $cAddr) {
            $handoffUrl = $hl.Address
            $handoffDisplay = $hl.TextToDisplay
        }
    }

    $handbackDate = $handbackInfo[$sheetName]

    foreach ($row in @(2, 3)) {
        # Status -> handed back, now in sync with en-US.
        $ws.Cells.Item($row, 2).Value = $newStatus

        # E: Latest Target File.
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $sourceUrl, "", "", $sourceDisplay) | Out-Null

        # F: Latest Handback File.
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $handoffUrl, "", "", $handoffDisplay) | Out-Null

        # G: Latest Handback DateTime.
        $ws.Cells.Item($row, 7).Value = $handbackDate
    }
}
